$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.090.37"
$ws.Range("E2").Value = "'  +7.69%  "
$ws.Range("D3").Value = "'1.874.02"
$ws.Range("E3").Value = "'  +5.44%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'249.03"
$ws.Range("E5").Value = "'  +2.60%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.01%  "
$ws.Range("D7").Value = "'0.4976"
$ws.Range("E7").Value = "'  +1.94%  "
$ws.Range("D8").Value = "'45.57"
$ws.Range("E8").Value = "'  +8.58%  "
$ws.Range("D9").Value = "'0.2843"
$ws.Range("E9").Value = "'  +7.32%  "
$ws.Range("D10").Value = "'0.06563"
$ws.Range("E10").Value = "'  +5.31%  "
$ws.Range("D11").Value = "'1.869.97"
$ws.Range("E11").Value = "'  +4.76%  "
$ws.Range("D12").Value = "'17.06"
$ws.Range("E12").Value = "'  +4.67%  "
$ws.Range("D13").Value = "'0.07190"
$ws.Range("E13").Value = "'  +2.66%  "
$ws.Range("D14").Value = "'0.6610"
$ws.Range("E14").Value = "'  +7.19%  "
$ws.Range("D15").Value = "'85.16"
$ws.Range("E15").Value = "'  +6.89%  "
$ws.Range("D16").Value = "'4.811"
$ws.Range("E16").Value = "'  +4.54%  "
$ws.Range("D17").Value = "'30.072.04"
$ws.Range("E17").Value = "'  +7.55%  "
$ws.Range("D18").Value = "'0.9995"
$ws.Range("E18").Value = "'  +0.03%  "
$ws.Range("D19").Value = "'12.91"
$ws.Range("E19").Value = "'  +9.38%  "
$ws.Range("D20").Value = "'0.000007514"
$ws.Range("E20").Value = "'  +4.37%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "'  -0.12%  "
$ws.Range("D22").Value = "'2.110.14"
$ws.Range("E22").Value = "'  +5.00%  "
$ws.Range("D23").Value = "'4.747"
$ws.Range("E23").Value = "'  +3.97%  "
$ws.Range("D24").Value = "'9.032"
$ws.Range("E24").Value = "'  +4.52%  "
$ws.Range("D25").Value = "'5.505"
$ws.Range("E25").Value = "'  +6.26%  "
$ws.Range("D26").Value = "'144.74"
$ws.Range("E26").Value = "'  +2.17%  "
$ws.Range("D27").Value = "'135.46"
$ws.Range("E27").Value = "'  +24.16%  "
$ws.Range("D28").Value = "'16.74"
$ws.Range("E28").Value = "'  +7.50%  "
$ws.Range("D29").Value = "'1.949"
$ws.Range("E29").Value = "'  +4.81%  "
$ws.Range("D30").Value = "'1.390"
$ws.Range("E30").Value = "'  -0.16%  "
$ws.Range("D31").Value = "'4.224"
$ws.Range("E31").Value = "'  +3.61%  "
$ws.Range("D32").Value = "'0.08603"
$ws.Range("E32").Value = "'  +4.27%  "
$ws.Range("D33").Value = "'3.894"
$ws.Range("E33").Value = "'  +3.46%  "
$ws.Range("D34").Value = "'0.05065"
$ws.Range("E34").Value = "'  +7.21%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "'  +6.95%  "
$ws.Range("D36").Value = "'0.6850"
$ws.Range("E36").Value = "'  +7.79%  "
$ws.Range("D37").Value = "'1.0000"
$ws.Range("D38").Value = "'2.698"
$ws.Range("E38").Value = "'  +3.63%  "
$ws.Range("D39").Value = "'2.335"
$ws.Range("E39").Value = "'  +14.16%  "
$ws.Range("D40").Value = "'2.740"
$ws.Range("E40").Value = "'  +5.64%  "
$ws.Range("D41").Value = "'0.9630"
$ws.Range("E41").Value = "'  +2.60%  "
$ws.Range("D42").Value = "'0.01627"
$ws.Range("E42").Value = "'  +6.40%  "
$ws.Range("D43").Value = "'6.077"
$ws.Range("E43").Value = "'  +3.92%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "'  +0.11%  "
$ws.Range("D45").Value = "'103.22"
$ws.Range("E45").Value = "'  +3.20%  "
$ws.Range("D46").Value = "'0.4182"
$ws.Range("E46").Value = "'  +6.45%  "
$ws.Range("D47").Value = "'7.498"
$ws.Range("E47").Value = "'  +5.02%  "
$ws.Range("D48").Value = "'0.1254"
$ws.Range("E48").Value = "'  +5.49%  "
$ws.Range("D49").Value = "'0.05636"
$ws.Range("E49").Value = "'  +4.22%  "
$ws.Range("D50").Value = "'32.53"
$ws.Range("E50").Value = "'  +7.13%  "
$ws.Range("D51").Value = "'8.268"
$ws.Range("E51").Value = "'  +3.80%  "
